$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1/Q1 with same style as the existing header cells ---
# Copy O1's format+value into P1:Q1 first so the new cells inherit the bold/border
# header style, then overwrite with the correct sequential values.
$ws.Range("O1").Copy($ws.Range("P1:Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    # Swap the I/K/M/O alternating 1/2 pattern
    $ws.Cells.Item($r, 9).Value = 2    # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1

    # New columns P and Q, unstyled, value 2
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
